# The commit swaps the deck's theme from the "Integral" (Red Violet) color
# scheme over to the stock "Office Theme" color scheme. The 13 slides (all
# built on the single SlideMaster / theme part the deck actually renders
# with) pick up the classic Office blue/orange palette in place of the
# magenta/purple Integral palette.
#
# Theme colours are addressed through the live ThemeColorScheme exposed on
# a slide (12 slots, in fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) and each slot's .RGB is assigned the packed R + G*256 + B*65536
# value for the target "Office" theme colour.

function PackRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (target), in clrScheme order.
$officeColors = @(
    (PackRGB 0x00 0x00 0x00),  # dk1      000000
    (PackRGB 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (PackRGB 0x44 0x54 0x6A),  # dk2      44546A
    (PackRGB 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (PackRGB 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (PackRGB 0xED 0x7D 0x31),  # accent2  ED7D31
    (PackRGB 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (PackRGB 0xFF 0xC0 0x00),  # accent4  FFC000
    (PackRGB 0x44 0x72 0xC4),  # accent5  4472C4
    (PackRGB 0x70 0xAD 0x47),  # accent6  70AD47
    (PackRGB 0x05 0x63 0xC1),  # hlink    0563C1
    (PackRGB 0x95 0x4F 0x72)   # folHlink 954F72
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
